$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting the existing
# "Late"/"Outstanding heading"/"Outstanding" columns one slot to the right.
$ws.Columns("N").Insert()

# Match the width Excel gives the freshly-inserted column (it copies the
# width of the column immediately to its left, column M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet (was "Transactions") and move
# the selection to M16.
$ws.Activate()
$ws.Range("M16").Select() | Out-Null
